$wb = $excel.ActiveWorkbook
$wsPresets = $wb.Worksheets.Item("Presets")

# Clear the "Testing Completed?" (X) marks in column B of the Presets sheet
# (the shared "X" string is then dropped automatically since it's no longer referenced,
# which shifts the later shared-string indices used on "Reference Lists" down by one)
$wsPresets.Range("B2:B9").ClearContents()

# Make "Presets" the active sheet/tab, with B2 as the selected cell
$wsPresets.Activate() | Out-Null
$wsPresets.Range("B2").Select() | Out-Null
